$wb = $excel.ActiveWorkbook

# Layer0 sheet
$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.5874250293444145
$ws0.Range("C2").Value = -1.166021715115316
$ws0.Range("B3").Value = -0.1806958593573982
$ws0.Range("C3").Value = 1.125158672491851
$ws0.Range("B4").Value = 0.9477571601906963
$ws0.Range("C4").Value = 0.9888840968501631

# Layer1 sheet
$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -1.711929502250102
$ws1.Range("C2").Value = 0.2787071566732836
$ws1.Range("B3").Value = 0.2854521736965125
$ws1.Range("C3").Value = -0.1904456298739174
$ws1.Range("B4").Value = 1.379515385541518
$ws1.Range("C4").Value = -0.413753189874506
